# Adds 4 new "Housekeeping" command-deck rows (45-48) to the all_systems
# sheet, for the formatter's new readout-intercept commands (0x30-0x34),
# by cloning the formatting of the last existing data row (44) and then
# filling in the new values/formula for each of the 4 new rows. Also
# clears the one-off AC-column styling that used to single out the last
# few "description" cells, since new rows keep that column unstyled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clone row 44's full formatting down into the 4 new rows ------------
$ws.Range("A44:AC44").Copy()
$ws.Range("A45:AC48").PasteSpecial(-4122)   # xlPasteFormats

# --- remove the one-off style override that used to live on AC34:AC44 ---
# (new rows leave AC unstyled, and the old distinguishing style goes away
# too)
$ws.Range("AC34:AC44").Style = "Normal"

# --- row 45: set_readout_disable (0x30) ----------------------------------
$ws.Range("A45").Value = "set_readout_disable"
$ws.Range("B45").Value = "0000 0000"
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 110000
$ws.Range("E45").Formula = '=_xlfn.CONCAT("0x", DEC2HEX(_xlfn.BITLSHIFT($C45,7) + BIN2DEC($D45)))'
$ws.Range("F45").Value = "—"
$ws.Range("G45").Value = "—"
$ws.Range("H45").Value = "—"
$ws.Range("I45").Value = "—"
$ws.Range("J45").Value = "—"
$ws.Range("K45").Value = "0x00"
$ws.Range("L45").Value = 1
$ws.Range("M45").Value = 0
$ws.Range("N45").Value = 0
$ws.Range("O45").Value = 0
$ws.Range("P45").Value = 0
$ws.Range("Q45").Value = 0
$ws.Range("R45").Value = 0
$ws.Range("S45").Value = 0
$ws.Range("T45").Value = 0
$ws.Range("U45").Value = 0
$ws.Range("V45").Value = 0
$ws.Range("W45").Value = 0
$ws.Range("X45").Value = "?"
$ws.Range("Y45").Value = "0x30"
$ws.Range("Z45").Value = "0x30"
$ws.Range("AA45").Value = "0x00"
$ws.Range("AB45").Value = "todo"
$ws.Range("AC45").Value = "Formatter must handle this, disabling HK all data readout"

# --- row 46: set_readout_enable_power (0x31) ------------------------------
$ws.Range("A46").Value = "set_readout_enable_power"
$ws.Range("B46").Value = "0000 0000"
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 110001
$ws.Range("E46").Formula = '=_xlfn.CONCAT("0x", DEC2HEX(_xlfn.BITLSHIFT($C46,7) + BIN2DEC($D46)))'
$ws.Range("F46").Value = "—"
$ws.Range("G46").Value = "—"
$ws.Range("H46").Value = "—"
$ws.Range("I46").Value = "—"
$ws.Range("J46").Value = "—"
$ws.Range("K46").Value = "0x00"
$ws.Range("L46").Value = 1
$ws.Range("M46").Value = 0
$ws.Range("N46").Value = 0
$ws.Range("O46").Value = 0
$ws.Range("P46").Value = 0
$ws.Range("Q46").Value = 0
$ws.Range("R46").Value = 0
$ws.Range("S46").Value = 0
$ws.Range("T46").Value = 0
$ws.Range("U46").Value = 0
$ws.Range("V46").Value = 0
$ws.Range("W46").Value = 0
$ws.Range("X46").Value = "?"
$ws.Range("Y46").Value = "0x30"
$ws.Range("Z46").Value = "0x31"
$ws.Range("AA46").Value = "0x00"
$ws.Range("AB46").Value = "todo"
$ws.Range("AC46").Value = "Formatter must handle this, enabling HK readout of AD7490"

# --- row 47: set_readout_enable_rtd (0x32) --------------------------------
$ws.Range("A47").Value = "set_readout_enable_rtd"
$ws.Range("B47").Value = "0000 0000"
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 110010
$ws.Range("E47").Formula = '=_xlfn.CONCAT("0x", DEC2HEX(_xlfn.BITLSHIFT($C47,7) + BIN2DEC($D47)))'
$ws.Range("F47").Value = "—"
$ws.Range("G47").Value = "—"
$ws.Range("H47").Value = "—"
$ws.Range("I47").Value = "—"
$ws.Range("J47").Value = "—"
$ws.Range("K47").Value = "0x00"
$ws.Range("L47").Value = 1
$ws.Range("M47").Value = 0
$ws.Range("N47").Value = 0
$ws.Range("O47").Value = 0
$ws.Range("P47").Value = 0
$ws.Range("Q47").Value = 0
$ws.Range("R47").Value = 0
$ws.Range("S47").Value = 0
$ws.Range("T47").Value = 0
$ws.Range("U47").Value = 0
$ws.Range("V47").Value = 0
$ws.Range("W47").Value = 0
$ws.Range("X47").Value = "?"
$ws.Range("Y47").Value = "0x30"
$ws.Range("Z47").Value = "0x32"
$ws.Range("AA47").Value = "0x00"
$ws.Range("AB47").Value = "todo"
$ws.Range("AC47").Value = "Formatter must handle this, enabling HK readout of RTDs"

# --- row 48: set_readout_enable_intro (0x34) ------------------------------
$ws.Range("A48").Value = "set_readout_enable_intro"
$ws.Range("B48").Value = "0000 0000"
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 110100
$ws.Range("E48").Formula = '=_xlfn.CONCAT("0x", DEC2HEX(_xlfn.BITLSHIFT($C48,7) + BIN2DEC($D48)))'
$ws.Range("F48").Value = "—"
$ws.Range("G48").Value = "—"
$ws.Range("H48").Value = "—"
$ws.Range("I48").Value = "—"
$ws.Range("J48").Value = "—"
$ws.Range("K48").Value = "0x00"
$ws.Range("L48").Value = 1
$ws.Range("M48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("O48").Value = 0
$ws.Range("P48").Value = 0
$ws.Range("Q48").Value = 0
$ws.Range("R48").Value = 0
$ws.Range("S48").Value = 0
$ws.Range("T48").Value = 0
$ws.Range("U48").Value = 0
$ws.Range("V48").Value = 0
$ws.Range("W48").Value = 0
$ws.Range("X48").Value = "?"
$ws.Range("Y48").Value = "0x30"
$ws.Range("Z48").Value = "0x34"
$ws.Range("AA48").Value = "0x00"
$ws.Range("AB48").Value = "todo"
$ws.Range("AC48").Value = "Formatter must handle this, enabling HK readout of introspection data"

# --- leave the selection where the author's did: AA48 --------------------
$ws.Range("AA48").Select()
